# LOB1040.docx: the "Objetivos"/"Docente(s)"/"Programa resumido"/"Programa"/
# "Avaliacao"/"Bibliografia" blocks got re-shuffled (content rotated between
# fixed heading slots). Every source string below is unique in the document,
# but several of the *new* values collide with *other* old values, so a
# straight sequential Find/Replace would cascade. We avoid that by first
# swapping every old value out to a unique placeholder token, then swapping
# each placeholder in for its final text. Scoping each Find to the specific
# paragraph keeps things extra safe.

$d = $word.ActiveDocument

function Replace-InParagraph($index, $old, $new) {
    $rng = $d.Paragraphs($index).Range
    $ok = $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Find failed for paragraph $index : $old"
    }
}

# ---- Step 1: move every affected value to a unique placeholder ----

Replace-InParagraph 6  "Capacitar o aluno no manuseio de medidores e circuitos de corrente contínua e alternada." "@@PH1@@"
Replace-InParagraph 7  "To enable the student in the circuits handling and DC/AC meters." "@@PH2@@"
Replace-InParagraph 9  "230696 - Carlos José Todero Peixoto" "@@PH3@@"
Replace-InParagraph 11 "Verificação experimental de aplicações em circuitos de corrente contínua e alternada." "@@PH4@@"
Replace-InParagraph 12 "Experimental verification of applications in DC and AC circuits." "@@PH5@@"
Replace-InParagraph 14 "1) Medidores. Osciloscópio.^l2) Tensão alternada.^l3) Potências. ^l4) Filtros.^l5) Ressonância.^l6) Campo magnético alternado." "@@PH6@@"
Replace-InParagraph 17 "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n." "@@PH7@@"
Replace-InParagraph 17 "NF≥ 5,0." "@@PH8@@"
Replace-InParagraph 17 "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada." "@@PH9@@"
Replace-InParagraph 19 "CAPUANO, G. Francisco; MARINO, M.A. Maria. Laboratório de eletricidade ^lEletrônica, Editora Érica (1998).^lMARKUS, Otávio. Circuitos elétricos: corrente contínua e corrente alternada-^lTeoria e Exercícios, Editora Érica, (2008).^lSADIKU, Mathew N. O.; ALEXANDER, Charles. Fundamentos de circuitos elétricos, Mcgraw-hill Interamericana (2009)." "@@PH10@@"

# ---- Step 2: swap each placeholder in for its final text ----

Replace-InParagraph 6  "@@PH1@@" "Verificação experimental de aplicações em circuitos de corrente contínua e alternada."
Replace-InParagraph 7  "@@PH2@@" "Experimental verification of applications in DC and AC circuits."
Replace-InParagraph 9  "@@PH3@@" "Capacitar o aluno no manuseio de medidores e circuitos de corrente contínua e alternada."
Replace-InParagraph 11 "@@PH4@@" "1) Medidores. Osciloscópio.^l2) Tensão alternada.^l3) Potências. ^l4) Filtros.^l5) Ressonância.^l6) Campo magnético alternado."
Replace-InParagraph 12 "@@PH5@@" "To enable the student in the circuits handling and DC/AC meters."
Replace-InParagraph 14 "@@PH6@@" "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
Replace-InParagraph 17 "@@PH7@@" "NF≥ 5,0."
Replace-InParagraph 17 "@@PH8@@" "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
Replace-InParagraph 17 "@@PH9@@" "CAPUANO, G. Francisco; MARINO, M.A. Maria. Laboratório de eletricidade ^lEletrônica, Editora Érica (1998).^lMARKUS, Otávio. Circuitos elétricos: corrente contínua e corrente alternada-^lTeoria e Exercícios, Editora Érica, (2008).^lSADIKU, Mathew N. O.; ALEXANDER, Charles. Fundamentos de circuitos elétricos, Mcgraw-hill Interamericana (2009)."
Replace-InParagraph 19 "@@PH10@@" "230696 - Carlos José Todero Peixoto"

Write-Output "Done"
